# Regenerate the localization-status report (CI "archive" snapshot).
#
# 1) Status moves from "Ready for handoff" -> "In Translation" everywhere it
#    appears: Overview!E2 (zh-cn column) + Overview!F2 (de-de column), and the
#    per-locale "Status" column (C2) on each of the "zh-cn" / "de-de" sheets.
#    Writing the same new text to every cell that held the old shared string
#    lets the workbook collapse back onto a single shared-string entry instead
#    of leaving the stale "Ready for handoff" string orphaned in the table.
#
# 2) Columns E/F on "Overview" and column C on "zh-cn"/"de-de" get narrower
#    (long-date columns were re-sized down from ~17.22 chars to ~13.41 chars
#    wide). Excel's ColumnWidth is in whole-character units of the Normal
#    style font (Calibri 11 here), so it only takes values on that font's
#    pixel grid -- we pick the grid point nearest the target width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# -- Status text: "Ready for handoff" -> "In Translation" -------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# -- Column widths: ~17.22 -> ~13.41 characters ------------------------------
$newColumnWidth = 12.5

$overview.Range("E1").EntireColumn.ColumnWidth = $newColumnWidth
$overview.Range("F1").EntireColumn.ColumnWidth = $newColumnWidth
$zhcn.Range("C1").EntireColumn.ColumnWidth = $newColumnWidth
$dede.Range("C1").EntireColumn.ColumnWidth = $newColumnWidth
